$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.123.90"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.655.28"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.97"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5261"
$ws.Range("E6").Value = "  +1.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06348"
$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.45"
$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07787"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.508"
$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").Value = "1.653.89"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5482"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("D15").Value = "0.0₅8226"
$ws.Range("E15").Value = "  +1.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.39"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "26.132.27"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.577"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.90"
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.08"
$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.043"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.93"
$ws.Range("E24").Value = "  +1.75%  "

$ws.Range("E25").Value = "  +2.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.269"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.16"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.433"
$ws.Range("E28").Value = "  +1.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05917"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.525"
$ws.Range("E31").Value = "  -1.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.255"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.585"
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9530"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.414"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5707"
$ws.Range("E37").Value = "  +1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01619"
$ws.Range("E38").Value = "  +2.13%  "

$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").Value = "1.031.22"
$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.80"
$ws.Range("E43").Value = "  +2.87%  "

$ws.Range("D44").Value = "1.801.89"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.34"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4300"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.479"
$ws.Range("E48").Value = "  +1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05163"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.829"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09728"
$ws.Range("E51").Value = "  +0.66%  "
